$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 844 (shifts old 844:865 down to 849:870)
$ws.Rows.Item(844).Resize(5).Insert()

# Constant columns for every row in this sheet (Mercado ID / Mercado / Region / Codreg / CategoriaID / Categoria / Clasificacion)
$constA = 10
$constB = "Vega Modelo de Temuco"
$constC = "La Araucanía"
$constE = 9
$constF = 100112020
$constG = "Tomate"
$constR = "Hortaliza"

function Set-TomateRow($Row, $D, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q) {
    $ws.Cells.Item($Row, 1).Value = $constA
    $ws.Cells.Item($Row, 2).Value = $constB
    $ws.Cells.Item($Row, 3).Value = $constC
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $constE
    $ws.Cells.Item($Row, 6).Value = $constF
    $ws.Cells.Item($Row, 7).Value = $constG
    $ws.Cells.Item($Row, 8).Value = $H
    $ws.Cells.Item($Row, 9).Value = $I
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = $constR
}

# New weekly rows for the 2021-09-09 (serial 44448) report, "Región de Arica y Parinacota"
Set-TomateRow 844 44448 "Larga vida" "Extra"   550  28000 30000 28909 "`$/bandeja 18 kilos" "Región de Arica y Parinacota" 1606 18
Set-TomateRow 845 44448 "Larga vida" "Extra"   200  32000 32000 32000 "`$/bandeja 20 kilos" "Región de Arica y Parinacota" 1600 20
Set-TomateRow 846 44448 "Larga vida" "Primera" 3500 24000 26000 24800 "`$/bandeja 18 kilos" "Región de Arica y Parinacota" 1378 18
Set-TomateRow 847 44448 "Larga vida" "Primera" 3000 10000 11000 10500 "`$/caja 10 kilos"    "Región de Arica y Parinacota" 1050 10
Set-TomateRow 848 44448 "Larga vida" "Segunda" 500  22000 22000 22000 "`$/bandeja 18 kilos" "Región de Arica y Parinacota" 1222 18
